$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 10 (shifts existing rows 10.. down by one)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the "Maquina de guerra" card
$ws.Cells.Item(10, 1).Value = "Maquina de guerra"
$ws.Cells.Item(10, 2).Value = "Federación de comercio"
$ws.Cells.Item(10, 3).Value = "Criatura"
$ws.Cells.Item(10, 4).Value = "Común"
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = 6
$ws.Cells.Item(10, 8).Value = "Entra en CB en PA cuando muere regresa en PD con 3/5"

# Comandante's ability text cell (now row 13 after the insert) becomes a bare number
$ws.Cells.Item(13, 8).Value = 66

# Match the saved selection from the authored workbook
$ws.Range("H11").Select()
